$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark completed to-do items with the green "done" fill (reuses existing style) ---
$ws.Range("B3").Interior.Color = 5287936   # win/lose action -> done (green FF00B050)
$ws.Range("B4").Interior.Color = 5287936   # bug fixing - load question -> done (green FF00B050)
$ws.Range("B5").Interior.Color = 5287936   # Try again on question -> done (green FF00B050)

# --- Highlight "friends" (E3) in light blue ---
$ws.Range("E3").Interior.Color = 15773696  # FF00B0F0

# --- Highlight "leaderboard" (D2) in gray (Text1 theme color) ---
$ws.Range("D2").Interior.Color = 5287936
$ws.Range("D2").Interior.ThemeColor = 1

# --- Add two new to-do items under "side menu" ---
$ws.Range("B6").Value = "log out button"
$ws.Range("B7").Value = "friends in the side menu"

# --- Grow Table1 to include the newly added row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E7"))

# --- Update the view: zoom level and active selection ---
$win = $wb.Windows.Item(1)
$win.Zoom = 130
$ws.Range("B8").Select()
